$wb = $excel.ActiveWorkbook

# ALC row 9
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 514.4286
$ws.Range("J9").Value = 92.5
$ws.Range("L9").Value = 92.5
$ws.Range("N9").Value = -430.5

# ALC row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 910.1667
$ws.Range("I98").Value = 910.1667
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 910.1667
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 587.8333
$ws.Range("N98").ClearContents()

# ALC row 103
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 974.75
$ws.Range("I103").Value = 949.5
$ws.Range("K103").Value = 2848.5
$ws.Range("M103").Value = -2262.5

# ALC row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 910.1667
$ws.Range("I122").Value = 910.1667
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2730.5001
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -280.5001000000002
$ws.Range("N122").ClearContents()

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1409.7646
$ws.Range("I132").Value = 1491
$ws.Range("J132").Value = 1030.6666
$ws.Range("K132").Value = 4473
$ws.Range("L132").Value = 3091.9998
$ws.Range("M132").Value = -1943
$ws.Range("N132").Value = -8151.9998

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 4017.7273
$ws.Range("I137").Value = 2112.25
$ws.Range("J137").Value = 5106.5713
$ws.Range("K137").Value = 6336.75
$ws.Range("L137").Value = 15319.7139
$ws.Range("M137").Value = -3786.75
$ws.Range("N137").Value = -20419.7139

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3080.46
$ws.Range("I138").Value = 1317.7727
$ws.Range("J138").Value = 4465.4287
$ws.Range("K138").Value = 3953.3181
$ws.Range("L138").Value = 13396.2861
$ws.Range("M138").Value = 1186.6819
$ws.Range("N138").Value = -23676.2861

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1137
$ws.Range("I2").Value = 705.5
$ws.Range("J2").Value = 2000
$ws.Range("K2").Value = 705.5
$ws.Range("L2").Value = 2000
$ws.Range("M2").Value = -592.5
$ws.Range("N2").Value = -2226

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1154.7142
$ws.Range("I74").Value = 1151.4166
$ws.Range("J74").Value = 1174.5
$ws.Range("K74").Value = 1151.4166
$ws.Range("L74").Value = 1174.5
$ws.Range("M74").Value = -277.4166
$ws.Range("N74").Value = -2922.5

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1154.7142
$ws.Range("I77").Value = 1151.4166
$ws.Range("J77").Value = 1174.5
$ws.Range("K77").Value = 5757.083000000001
$ws.Range("L77").Value = 5872.5
$ws.Range("M77").Value = -1389.083000000001
$ws.Range("N77").Value = -14608.5

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1137
$ws.Range("I116").Value = 705.5
$ws.Range("J116").Value = 2000
$ws.Range("K116").Value = 705.5
$ws.Range("L116").Value = 2000
$ws.Range("M116").Value = 1588.5
$ws.Range("N116").Value = -6588

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1137
$ws.Range("I3").Value = 705.5
$ws.Range("J3").Value = 2000
$ws.Range("K3").Value = 705.5
$ws.Range("L3").Value = 2000
$ws.Range("M3").Value = -591.5
$ws.Range("N3").Value = -2228

# BSM row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3667.6875
$ws.Range("I20").Value = 3880.7273
$ws.Range("K20").Value = 3880.7273
$ws.Range("M20").Value = -3633.7273

# BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 495.5
$ws.Range("I94").Value = 495.5
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 495.5
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -44.5
$ws.Range("N94").ClearContents()

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2688.9285
$ws.Range("I105").Value = 2665
$ws.Range("K105").Value = 2665
$ws.Range("M105").Value = -918

# CRP row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 992.3333
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 992.3333
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 992.3333
$ws.Range("N22").Value = -1692.3333
$ws.Range("M22").ClearContents()

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2061.5454
$ws.Range("I31").Value = 2118.1
$ws.Range("K31").Value = 2118.1
$ws.Range("M31").Value = -1823.1

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2061.5454
$ws.Range("I34").Value = 2118.1
$ws.Range("K34").Value = 2118.1
$ws.Range("M34").Value = -1916.1

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1948.7667
$ws.Range("I58").Value = 1933.5
$ws.Range("J58").Value = 1990.75
$ws.Range("K58").Value = 1933.5
$ws.Range("L58").Value = 1990.75
$ws.Range("M58").Value = -1730.5
$ws.Range("N58").Value = -2396.75

# CRP row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 2148.6
$ws.Range("I107").Value = 1447
$ws.Range("K107").Value = 1447
$ws.Range("M107").Value = 473

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1484.6842
$ws.Range("I134").Value = 1522.7778
$ws.Range("J134").Value = 799
$ws.Range("K134").Value = 4568.3334
$ws.Range("L134").Value = 2397
$ws.Range("M134").Value = -2033.3334
$ws.Range("N134").Value = -7467

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1948.7667
$ws.Range("I136").Value = 1933.5
$ws.Range("J136").Value = 1990.75
$ws.Range("K136").Value = 5800.5
$ws.Range("L136").Value = 5972.25
$ws.Range("M136").Value = -3250.5
$ws.Range("N136").Value = -11072.25

# CUL row 86
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 294.66666
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

# CUL row 89
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H89").Value = 294.66666
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

# CUL row 107
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1162.091
$ws.Range("J107").Value = 278.3
$ws.Range("L107").Value = 834.9000000000001
$ws.Range("N107").Value = -4674.9

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1220.1666
$ws.Range("J131").Value = 1749.5
$ws.Range("L131").Value = 5248.5
$ws.Range("N131").Value = -15328.5

# CUL row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 4234.875
$ws.Range("I132").Value = 3776.6
$ws.Range("J132").Value = 4998.6665
$ws.Range("K132").Value = 33989.4
$ws.Range("L132").Value = 44987.9985
$ws.Range("M132").Value = -31459.4
$ws.Range("N132").Value = -50047.9985

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7844.4614
$ws.Range("I70").Value = 8227.9
$ws.Range("K70").Value = 8227.9
$ws.Range("M70").Value = -7957.9

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 7844.4614
$ws.Range("I73").Value = 8227.9
$ws.Range("K73").Value = 8227.9
$ws.Range("M73").Value = -7291.9

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3642.2222
$ws.Range("I102").Value = 3347.625
$ws.Range("J102").Value = 5999
$ws.Range("K102").Value = 3347.625
$ws.Range("L102").Value = 5999
$ws.Range("M102").Value = -1725.625
$ws.Range("N102").Value = -9243

# LTW row 16
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 667
$ws.Range("I16").Value = 583.75
$ws.Range("K16").Value = 583.75
$ws.Range("M16").Value = -413.75

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2849.9443
$ws.Range("I132").Value = 1572.4286
$ws.Range("J132").Value = 3662.9092
$ws.Range("K132").Value = 4717.2858
$ws.Range("L132").Value = 10988.7276
$ws.Range("M132").Value = -2187.2858
$ws.Range("N132").Value = -16048.7276

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2668.875
$ws.Range("I126").Value = 2742.8333
$ws.Range("J126").Value = 2447
$ws.Range("K126").Value = 8228.499899999999
$ws.Range("L126").Value = 7341
$ws.Range("M126").Value = -5758.499899999999
$ws.Range("N126").Value = -12281
